# Apply updated cryptocurrency price/volume data per the upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.258.28"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "3.493.18"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.13"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.09"
$ws.Range("E6").Value = "  +1.77%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.80%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("E10").Value = "  +1.87%  "

$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("D12").Value = "4.087.15"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("E13").Value = "  +1.91%  "

$ws.Range("E14").Value = "  +1.07%  "

$ws.Range("D15").Value = "3.492.04"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.73"
$ws.Range("E16").Value = "  -6.65%  "

$ws.Range("D17").Value = "64.334.17"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.87"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.74"
$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.63"
$ws.Range("E21").Value = "  +2.78%  "

$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").Value = "3.633.25"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.29"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000114"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.50"
$ws.Range("E29").Value = "  -4.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("D32").Value = "3.513.36"
$ws.Range("E32").Value = "  -0.26%  "

$ws.Range("E33").Value = "  +3.64%  "

$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.39"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.15"
$ws.Range("E36").Value = "  -3.75%  "

$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.93"
$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0781"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("E41").Value = "  -0.96%  "

$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  -4.75%  "

$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("E45").Value = "  +2.03%  "

$ws.Range("E46").Value = "  -3.38%  "

$ws.Range("D47").Value = "2.459.70"
$ws.Range("E47").Value = "  +0.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.897"
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("E50").Value = "  -1.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.14"
$ws.Range("E51").Value = "  -1.40%  "
